$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 26, shifting rows 26:35 down to 27:36,
# then fill in the new row's values (a new weekly "Especial" price entry).
$ws.Rows("26:26").Insert(-4121)  # -4121 = xlShiftDown

$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44985
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100107
$ws.Cells.Item(26, 8).Value = "Otros"
$ws.Cells.Item(26, 9).Value = 100107011
$ws.Cells.Item(26, 10).Value = "Tuna"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Especial"
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 13000
$ws.Cells.Item(26, 15).Value = 14000
$ws.Cells.Item(26, 16).Value = 13500
$ws.Cells.Item(26, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(26, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 19).Value = 750
$ws.Cells.Item(26, 20).Value = 18
